$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4677.5
$ws.Range("I2").Value = 3608.3333
$ws.Range("J2").Value = 6815.8335
$ws.Range("K2").Value = 3608.3333
$ws.Range("L2").Value = 6815.8335
$ws.Range("M2").Value = -3495.3333
$ws.Range("N2").Value = -7041.8335
$ws.Range("H11").Value = 10209.75
$ws.Range("I11").Value = 10209.75
$ws.Range("K11").Value = 10209.75
$ws.Range("M11").Value = -10069.75
$ws.Range("H43").Value = 144446110
$ws.Range("I43").Value = 144446110
$ws.Range("K43").Value = 144446110
$ws.Range("M43").Value = -144446041
$ws.Range("H46").Value = 4333
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2881
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("K47").Value = 20000
$ws.Range("M47").Value = -19028
$ws.Range("H60").Value = 4333
$ws.Range("I60").Value = 1000
$ws.Range("K60").Value = 3000
$ws.Range("M60").Value = -2516
$ws.Range("H62").Value = 1703.8889
$ws.Range("J62").Value = 2523.5
$ws.Range("L62").Value = 2523.5
$ws.Range("N62").Value = -3771.5
$ws.Range("H65").Value = 1703.8889
$ws.Range("J65").Value = 2523.5
$ws.Range("L65").Value = 12617.5
$ws.Range("N65").Value = -18857.5
$ws.Range("H69").Value = 13139.77
$ws.Range("I69").Value = 9255.571
$ws.Range("K69").Value = 27766.713
$ws.Range("M69").Value = -26892.713
$ws.Range("H72").Value = 13139.77
$ws.Range("I72").Value = 9255.571
$ws.Range("K72").Value = 83300.139
$ws.Range("M72").Value = -78932.139
$ws.Range("H88").Value = 1693
$ws.Range("J88").Value = 1693
$ws.Range("L88").Value = 1693
$ws.Range("N88").Value = -2505
$ws.Range("H91").Value = 1693
$ws.Range("J91").Value = 1693
$ws.Range("L91").Value = 1693
$ws.Range("N91").Value = -4501
$ws.Range("H103").Value = 480.75
$ws.Range("I103").Value = 313.33334
$ws.Range("J103").Value = 648.1667
$ws.Range("K103").Value = 940.0000200000001
$ws.Range("L103").Value = 1944.5001
$ws.Range("M103").Value = -354.0000200000001
$ws.Range("N103").Value = -3116.5001
$ws.Range("H104").Value = 698.6
$ws.Range("I104").Value = 373.25
$ws.Range("J104").Value = 2000
$ws.Range("K104").Value = 1119.75
$ws.Range("L104").Value = 6000
$ws.Range("M104").Value = 627.25
$ws.Range("N104").Value = -9494
$ws.Range("H116").Value = 5031.154
$ws.Range("I116").Value = 5166.4
$ws.Range("K116").Value = 5166.4
$ws.Range("M116").Value = -1724.4
$ws.Range("H137").Value = 3854424
$ws.Range("I137").Value = 7019.4
$ws.Range("J137").Value = 9100885
$ws.Range("K137").Value = 21058.2
$ws.Range("L137").Value = 27302655
$ws.Range("M137").Value = -18508.2
$ws.Range("N137").Value = -27307755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7158.8
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7158.8
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7158.8
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -7796.8
$ws.Range("H74").Value = 587691.1
$ws.Range("I74").Value = 1062.8889
$ws.Range("K74").Value = 1062.8889
$ws.Range("M74").Value = -188.8888999999999
$ws.Range("H77").Value = 587691.1
$ws.Range("I77").Value = 1062.8889
$ws.Range("K77").Value = 5314.4445
$ws.Range("M77").Value = -946.4444999999996
$ws.Range("H86").Value = 144166.67
$ws.Range("J86").Value = 144166.67
$ws.Range("L86").Value = 144166.67
$ws.Range("N86").Value = -146538.67
$ws.Range("H89").Value = 144166.67
$ws.Range("J89").Value = 144166.67
$ws.Range("L89").Value = 432500.01
$ws.Range("N89").Value = -444356.01
$ws.Range("H102").Value = 2211.9048
$ws.Range("I102").Value = 2211.9048
$ws.Range("K102").Value = 2211.9048
$ws.Range("M102").Value = -589.9047999999998
$ws.Range("H122").Value = 1966.5
$ws.Range("I122").Value = 1933
$ws.Range("K122").Value = 5799
$ws.Range("M122").Value = -3349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 872.3684
$ws.Range("J20").Value = 850.9
$ws.Range("L20").Value = 850.9
$ws.Range("N20").Value = -1344.9
$ws.Range("H80").Value = 166670830
$ws.Range("J80").Value = 4996.5
$ws.Range("L80").Value = 4996.5
$ws.Range("N80").Value = -6992.5
$ws.Range("H83").Value = 166670830
$ws.Range("J83").Value = 4996.5
$ws.Range("L83").Value = 24982.5
$ws.Range("N83").Value = -34966.5
$ws.Range("H86").Value = 5755.5454
$ws.Range("I86").Value = 3600.625
$ws.Range("K86").Value = 3600.625
$ws.Range("M86").Value = -2477.625
$ws.Range("H89").Value = 5755.5454
$ws.Range("I89").Value = 3600.625
$ws.Range("K89").Value = 18003.125
$ws.Range("M89").Value = -12387.125
$ws.Range("H105").Value = 39926
$ws.Range("I105").Value = 50389.5
$ws.Range("J105").Value = 18999
$ws.Range("K105").Value = 50389.5
$ws.Range("L105").Value = 18999
$ws.Range("M105").Value = -48642.5
$ws.Range("N105").Value = -22493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H17").Value = 9247
$ws.Range("I17").Value = 9247
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 9247
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -9073
$ws.Range("N17").ClearContents()
$ws.Range("H21").Value = 29790
$ws.Range("J21").Value = 29790
$ws.Range("L21").Value = 29790
$ws.Range("N21").Value = -30260
$ws.Range("H37").Value = 8450
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H107").Value = 1010.5172
$ws.Range("J107").Value = 627.5714
$ws.Range("L107").Value = 627.5714
$ws.Range("N107").Value = -4467.5714
$ws.Range("H122").Value = 3570.64
$ws.Range("I122").Value = 3420.5
$ws.Range("J122").Value = 3956.7144
$ws.Range("K122").Value = 10261.5
$ws.Range("L122").Value = 11870.1432
$ws.Range("M122").Value = -7811.5
$ws.Range("N122").Value = -16770.1432
$ws.Range("H132").Value = 13934201
$ws.Range("J132").Value = 37040270
$ws.Range("L132").Value = 111120810
$ws.Range("N132").Value = -111125870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8552552
$ws.Range("I56").Value = 8552552
$ws.Range("K56").Value = 8552552
$ws.Range("M56").Value = -8552022

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6000.4707
$ws.Range("J2").Value = 10134.9
$ws.Range("L2").Value = 10134.9
$ws.Range("N2").Value = -10360.9
$ws.Range("H33").Value = 18
$ws.Range("I33").Value = 18
$ws.Range("K33").Value = 18
$ws.Range("M33").Value = 234
$ws.Range("H102").Value = 26316740
$ws.Range("I102").Value = 26316740
$ws.Range("K102").Value = 26316740
$ws.Range("M102").Value = -26315118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4831.95
$ws.Range("I22").Value = 1329.125
$ws.Range("J22").Value = 7167.1665
$ws.Range("K22").Value = 1329.125
$ws.Range("L22").Value = 7167.1665
$ws.Range("M22").Value = -1034.125
$ws.Range("N22").Value = -7757.1665
$ws.Range("H27").Value = 4831.95
$ws.Range("I27").Value = 1329.125
$ws.Range("J27").Value = 7167.1665
$ws.Range("K27").Value = 1329.125
$ws.Range("L27").Value = 7167.1665
$ws.Range("M27").Value = -1222.125
$ws.Range("N27").Value = -7381.1665
$ws.Range("H46").Value = 14207.833
$ws.Range("J46").Value = 6666.4443
$ws.Range("L46").Value = 6666.4443
$ws.Range("N46").Value = -7042.4443
$ws.Range("H68").Value = 3777.7144
$ws.Range("I68").Value = 2666.3333
$ws.Range("K68").Value = 2666.3333
$ws.Range("M68").Value = -1917.3333
$ws.Range("H71").Value = 3777.7144
$ws.Range("I71").Value = 2666.3333
$ws.Range("K71").Value = 13331.6665
$ws.Range("M71").Value = -9587.666499999999
$ws.Range("H93").Value = 1436
$ws.Range("I93").Value = 1244.2858
$ws.Range("K93").Value = 1244.2858
$ws.Range("M93").Value = 3.714199999999892
$ws.Range("H100").Value = 3298.3333
$ws.Range("I100").Value = 2732.25
$ws.Range("J100").Value = 3751.2
$ws.Range("K100").Value = 2732.25
$ws.Range("L100").Value = 3751.2
$ws.Range("M100").Value = -2191.25
$ws.Range("N100").Value = -4833.2
$ws.Range("H132").Value = 3595.5334
$ws.Range("I132").Value = 3507.875
$ws.Range("K132").Value = 10523.625
$ws.Range("M132").Value = -7993.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 7143254.5
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 14286214
$ws.Range("K107").Value = 885
$ws.Range("L107").Value = 42858642
$ws.Range("M107").Value = 1035
$ws.Range("N107").Value = -42862482
$ws.Range("H133").Value = 37600
$ws.Range("J133").Value = 37600
$ws.Range("L133").Value = 37600
$ws.Range("N133").Value = -47720
